# Append newly-scraped Lancers listings (2025-09-22 01:22 JST run) and
# re-sort by priority score, matching the upstream scraper's behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-22 01:22:01"

# --- Step 1: insert one row before the old row 5 (GAS) for the new
#     "自己分析アプリ" listing (score 93, slots in after 検査報告書/123). ---
$ws.Rows("5:5").Insert()

$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【急募】自己分析アプリのバックエンド開発アドバイザリー募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5397930"
$ws.Range("G5").Value = 93
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# --- Step 2: insert two rows before the old row 8 (now row 9, Gemini)
#     for the new "データセンター" (score 18) and "中小企業支援"
#     (score 10) listings. ---
$ws.Rows("9:10").Insert()

$ws.Range("A9").Value = $newTimestamp
$ws.Range("B9").Value = "データセンター向けサーバー・ルーター設置作業"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5397887"
$ws.Range("G9").Value = 18

$ws.Range("A10").Value = $newTimestamp
$ws.Range("B10").Value = "【中小企業支援】債務超過・赤字経営解消の診断依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5397962"
$ws.Range("G10").Value = 10

# --- Refresh the "取得日時" timestamp on every row (new + pre-existing). ---
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp
$ws.Range("A11").Value = $newTimestamp
$ws.Range("A12").Value = $newTimestamp

# --- Register hyperlinks only for the cells beyond the worksheet's
#     original hyperlink coverage (F2:F9 already carry relationship ids
#     rId1-rId8, which the row-insert operations leave pointing at their
#     original - now displaced - rows, exactly as upstream produced it).
#     Restore the standard "URL column" formatting afterwards, since
#     Hyperlinks.Add otherwise stamps its own style. ---
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5397962")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5397817")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5395809")

$ws.Range("F2").Copy()
$ws.Range("F10:F12").PasteSpecial(-4122)
$excel.CutCopyMode = 0
